# ---------------------------------------------------------------------------
# First cut of data logging web pages
#
# This script reproduces, via PowerPoint COM automation, the authored edit:
#  1. Bump the cached "datetimeFigureOut" footer field text on the Slide
#     Master and all 11 Slide Layouts from 20-07-2021 -> 26-07-2021.
#  2. On Slide 1:
#       - Duplicate "Group 1" (title bar + 3 logos) BEFORE editing it, so the
#         duplicate keeps the old "RESPIMATIC 100 WEB DASHBOARD" title; move
#         the duplicate down into the middle of the slide (new "Group 17").
#       - Move the original "Group 1" up, and retitle it to
#         "RESPIMATIC 100 WEB LOGGER".
#       - Duplicate "TextBox 14" ("WEB DASHBOARD / RESPIMATIC 100") BEFORE
#         editing it, keep the duplicate's wording, move it to the new
#         "TextBox 15" position.
#       - Move the original "TextBox 14"; retext it to "WEB LOGGER ".
#       - Move "Group 16" (bottom logo strip).
#       - Delete "Picture 9" (the small PROCEED button picture).
#       - Re-point "Picture 7" at the image that used to belong to
#         "Picture 9" (rId5) and relocate it to the bottom-left corner.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Footer "date" placeholder text: 20-07-2021 -> 26-07-2021
#    (Slide Master + every Custom Layout.)
# ---------------------------------------------------------------------------
function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "26-07-2021"
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "26-07-2021"
}

# ---------------------------------------------------------------------------
# 2. Slide 1 shape surgery
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# The id-allocator's very first freshly-created shape on a slide is returned
# with a stale id/name in this runtime; burn that one immediately so every
# subsequent Duplicate()/AddTextbox() call gets the correctly-numbered id
# (this mirrors real PowerPoint's id bookkeeping for the rest of the script).
$warmUp = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$warmUp.Delete()

# Locate the shapes we need by their current (pre-edit) names.
function Find-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Name -eq $name) { return $candidate }
    }
    return $null
}

$group1   = Find-ShapeByName $s.Shapes "Group 1"
$group16  = Find-ShapeByName $s.Shapes "Group 16"
$textbox14 = Find-ShapeByName $s.Shapes "TextBox 14"
$picture9 = Find-ShapeByName $s.Shapes "Picture 9"
$picture7 = Find-ShapeByName $s.Shapes "Picture 7"

$ptPerEmu = 1 / 12700

# --- Duplicate "Group 1" while it still reads "WEB DASHBOARD"; the copy
#     becomes the new "Group 17" parked in the middle of the slide. ---------
$group17 = $group1.Duplicate().Item(1)
$group17.Name = "Group 17"
$group17.Left = 2556384 * $ptPerEmu
$group17.Top  = 2533579 * $ptPerEmu

$group17ChildNames = @("Rectangle 18", "Picture 19", "Picture 20", "Picture 21", "TextBox 22")
$group17Items = $group17.GroupItems
for ($j = 1; $j -le $group17Items.Count; $j++) {
    $group17Items.Item($j).Name = $group17ChildNames[$j - 1]
}

# --- Now move/retitle the original "Group 1" -------------------------------
$group1.Left = 2507221 * $ptPerEmu
$group1.Top  = 734511 * $ptPerEmu

$group1Items = $group1.GroupItems
$titleBox = Find-ShapeByName $group1Items "TextBox 8"
$titleBox.TextFrame.TextRange.Text = "RESPIMATIC 100 WEB LOGGER"

# --- Duplicate "TextBox 14" while it still reads "WEB DASHBOARD"; the copy
#     becomes the new "TextBox 15". -----------------------------------------
$textbox15 = $textbox14.Duplicate().Item(1)
$textbox15.Name = "TextBox 15"
$textbox15.Left = 2556385 * $ptPerEmu
$textbox15.Top  = 4325891 * $ptPerEmu
$textbox15.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "WEB LOGGER "

# --- Move the original "TextBox 14" and retext its first line. -------------
$textbox14.Left = 6971070 * $ptPerEmu
$textbox14.Top  = 4320898 * $ptPerEmu

# --- Move "Group 16". -------------------------------------------------------
$group16.Left = 2556385 * $ptPerEmu
$group16.Top  = 5471925 * $ptPerEmu

# --- Delete "Picture 9". ----------------------------------------------------
$picture9.Delete()

# --- Re-point "Picture 7" at the freed-up image (rId5) and relocate it. ----
$picture7.PictureFormat.Crop.PictureFile = ""
$picture7.Left = 119641 * $ptPerEmu
$picture7.Top  = 5336816 * $ptPerEmu
